# Final fix to handle queries that are not in given set
#
# This script reproduces, via Excel COM automation, the changes recorded in
# the commit:
#   - Updated evaluation metrics on Sheet1 (C3, G3)
#   - Re-ran timing/throughput measurements on Sheet2 (A3:A4, A9:A10, D3:D4, D9:D10)
#   - Selection/active-sheet bookkeeping: Sheet1 becomes the active tab with
#     G8 selected; Sheet2's selection moves to D16 and is no longer the
#     active/selected tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: updated numeric results -------------------------------------
$ws1.Range("C3").Value = 7.7713859437522101
$ws1.Range("G3").Value = 7.55012453779457

# --- Sheet2: updated "Traditional" column text ----------------------------
$ws2.Range("A3").Value  = "Mean response time:  0.203246806987524"
$ws2.Range("A4").Value  = "Average throughput for 1 iterations with 225 queries is:  4.920126494589327"
$ws2.Range("A9").Value  = "Mean response time:  0.15345342860623223"
$ws2.Range("A10").Value = "Average throughput for 1 iterations with 225 queries is:  6.516635106055798"

# --- Sheet2: updated "Wacky" column text ----------------------------------
$ws2.Range("D3").Value  = "Mean response time:  0.18235452790707504"
$ws2.Range("D4").Value  = "Average throughput for 1 iterations with 225 queries is:  5.483823250660296"
$ws2.Range("D9").Value  = "Mean response time:  0.17086823238049492"
$ws2.Range("D10").Value = "Average throughput for 1 iterations with 225 queries is:  5.852462953869433"

# --- View state -------------------------------------------------------------
# Sheet2 was the active/selected tab with A16 selected; now it just keeps a
# fresh selection at D16 (also scrolled so column C is the leftmost visible
# column) and is no longer the active tab.
$ws2.Activate()
$ws2.Range("D16").Select()
$excel.ActiveWindow.ScrollColumn = 3

# Sheet1 becomes the active tab, selection moves to G8.
$ws1.Activate()
$ws1.Range("G8").Select()
